# Applies the CaseJourneyDetails.xlsx update:
#   - Case Id (column A) changes from CASE329 to CASE383 for rows 2-5
#   - Received/Target/Creation/Start/End timestamps (columns G,H,I,N,O)
#     shift from Sep/2024 to Oct/2024 for rows 2-5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "CASE383"
$ws.Range("G2").Value = "22-Oct-2024 09:06:52 AM"
$ws.Range("H2").Value = "29-Oct-2024 09:06:52 AM"
$ws.Range("I2").Value = "22-Oct-2024 09:06:53 AM"
$ws.Range("N2").Value = "22-Oct-2024 09:07:39 AM"
$ws.Range("O2").Value = "22-Oct-2024 09:07:41 AM"

# Row 3
$ws.Range("A3").Value = "CASE383"
$ws.Range("G3").Value = "22-Oct-2024 09:06:52 AM"
$ws.Range("H3").Value = "29-Oct-2024 09:06:52 AM"
$ws.Range("I3").Value = "22-Oct-2024 09:06:53 AM"
$ws.Range("N3").Value = "22-Oct-2024 09:07:25 AM"
$ws.Range("O3").Value = "22-Oct-2024 09:07:25 AM"

# Row 4
$ws.Range("A4").Value = "CASE383"
$ws.Range("G4").Value = "22-Oct-2024 09:06:52 AM"
$ws.Range("H4").Value = "29-Oct-2024 09:06:52 AM"
$ws.Range("I4").Value = "22-Oct-2024 09:06:53 AM"
$ws.Range("N4").Value = "22-Oct-2024 09:07:04 AM"
$ws.Range("O4").Value = "22-Oct-2024 09:07:09 AM"

# Row 5
$ws.Range("A5").Value = "CASE383"
$ws.Range("G5").Value = "22-Oct-2024 09:06:52 AM"
$ws.Range("H5").Value = "29-Oct-2024 09:06:52 AM"
$ws.Range("I5").Value = "22-Oct-2024 09:06:53 AM"
$ws.Range("N5").Value = "22-Oct-2024 09:06:53 AM"
$ws.Range("O5").Value = "22-Oct-2024 09:06:53 AM"
